$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Tarantula")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 6.523762226400829
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 21.75520932085312

$ws = $wb.Worksheets.Item("Ochiai")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 3.685184637569725
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 19.82851293219664

$ws = $wb.Worksheets.Item("Op2")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 5.327431842653215
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 21.29711610884301

$ws = $wb.Worksheets.Item("Barinel")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 6.530633189709885
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 21.75864480250765

$ws = $wb.Worksheets.Item("Dstar")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 3.760243388654732
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 20.23220377015845

$ws = $wb.Worksheets.Item("Russell_rao")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 10.67334570686556
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 26.99384048580321

$ws = $wb.Worksheets.Item("Simple_matching")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 20.31487276541489
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 54.23551400893749

$ws = $wb.Worksheets.Item("Rogers_tanimoto")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 20.31487276541489
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 54.23551400893749

$ws = $wb.Worksheets.Item("Ample")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 3.800077581003437
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 28.61773613065617

$ws = $wb.Worksheets.Item("Jaccard")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 4.041909397216465
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 19.13158851453251

$ws = $wb.Worksheets.Item("Cohen")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 4.186112652487536
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 19.96867188627947

$ws = $wb.Worksheets.Item("Scott")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 5.743342558398804
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 20.83463069876829

$ws = $wb.Worksheets.Item("Rogot1")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 5.743342558398804
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 20.83463069876829

$ws = $wb.Worksheets.Item("Geometric_mean")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 3.725236265466182
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 20.33644237175215

$ws = $wb.Worksheets.Item("M2")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 4.406070452596446
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 21.20927214755002

$ws = $wb.Worksheets.Item("Wong1")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 10.67334570686556
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 26.99384048580321

$ws = $wb.Worksheets.Item("Sokal")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 20.31487276541489
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 54.23551400893749

$ws = $wb.Worksheets.Item("Sorensen_dice")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 4.041909397216465
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 19.13158851453251

$ws = $wb.Worksheets.Item("Dice")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 4.041909397216465
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 19.13158851453251

$ws = $wb.Worksheets.Item("Humman")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 20.31487276541489
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 54.23551400893749

$ws = $wb.Worksheets.Item("Wong2")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 20.31487276541489
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 54.23551400893749

$ws = $wb.Worksheets.Item("Euclid")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 20.31487276541489
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 54.23551400893749

$ws = $wb.Worksheets.Item("Zoltar")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 5.266201992405422
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 20.93808653239007

$ws = $wb.Worksheets.Item("Rogot2")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 3.909404174414617
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 20.32374413576328

$ws = $wb.Worksheets.Item("Hamming")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 20.31487276541489
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 54.23551400893749

$ws = $wb.Worksheets.Item("Fleiss")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 5.52808136612143
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 21.68362953853219

$ws = $wb.Worksheets.Item("Anderberg")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 4.041909397216465
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 19.13158851453251

$ws = $wb.Worksheets.Item("Goodman")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 4.041909397216465
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 19.13158851453251

$ws = $wb.Worksheets.Item("Harmonic_mean")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 3.909404174414617
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 20.58092690164781

$ws = $wb.Worksheets.Item("Kulczynski2")
$ws.Cells.Item(14, 1).Value = "Best exam"
$ws.Cells.Item(14, 3).Value = 4.246646708982527
$ws.Cells.Item(15, 1).Value = "Worst exam"
$ws.Cells.Item(15, 3).Value = 21.21749121125517

